# v0.4.3 BOM update: minor parts updates
#  - Row 17: MOV-20D220K Surge Absorber -> V14E14AUTO Varistor (Bournes)
#  - Row 33: 160R resistor -> 10R resistor (and drop the Mouser note/part)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17: replace surge-absorber part with the new varistor part ---
$ws.Range("F17").Value = "VARISTOR 22V 3KA DISC 14MM"
$ws.Range("L17").Value = "V14E14AUTO"
$ws.Range("M17").Value = "F6389-ND"

# --- Row 33: replace 160R resistor with the new 10R resistor ---
$ws.Range("F33").Value = "RES 10 OHM 2W 5% AXIAL"
$ws.Range("L33").Value = "FMP200JR-52-10R"
$ws.Range("M33").Value = "10ZCT-ND"
$ws.Range("N33").Value = ""

# --- Hyperlinks: the old M17 (surge absorber) part had a Digikey hyperlink
#     that doesn't apply to the new part, so it must be dropped. The engine's
#     Hyperlinks collection only supports a sheet-wide Delete(), so capture
#     the remaining links first, wipe them all, then recreate everybody
#     except the M17 one. ---
$remaining = @(
  @{Cell="M6";  Addr="http://search.digikey.com/us/en/products/TAP476K010SCS/478-1910-ND/564013"; Disp="478-1910-ND"},
  @{Cell="M13"; Addr="http://search.digikey.com/us/en/products/1N5818-TP/1N5818-TPCT-ND/950587";   Disp=""},
  @{Cell="M26"; Addr="http://search.digikey.com/us/en/products/MFR-25FBF-10K0/10.0KXBK-ND/13219";   Disp=""},
  @{Cell="M30"; Addr="http://search.digikey.com/us/en/products/RC55Y-2K49BI/985-1047-1-ND/2401912"; Disp="985-1047-1-ND"},
  @{Cell="M31"; Addr="http://search.digikey.com/us/en/products/MFP-25BRD52-3K9/3.9KADCT-ND/2059137";Disp=""},
  @{Cell="M36"; Addr="http://www.digikey.com.au/product-detail/en/MPX4250AP/MPX4250AP-ND/464053";   Disp=""},
  @{Cell="M3";  Addr="http://www.digikey.com/product-detail/en/TAP106K035SCS/478-1842-ND/563945";   Disp="478-1842-ND"},
  @{Cell="M7";  Addr="http://search.digikey.com/us/en/products/FK14X7R1H334K/445-5312-ND/2256792";  Disp="445-5312-ND"},
  @{Cell="M8";  Addr="http://www.digikey.com/product-detail/en/C315C103K5R5TA/399-4148-ND/817924";  Disp="399-4148-ND"},
  @{Cell="M32"; Addr="http://search.digikey.com/us/en/products/MFR-25FBF-100K/100KXBK-ND/13473";    Disp=""}
)

$ws.Hyperlinks.Delete()

foreach ($link in $remaining) {
    $r = $ws.Range($link.Cell)
    if ($link.Disp -ne "") {
        $r.Hyperlinks.Add($r, $link.Addr, "", "", $link.Disp)
    } else {
        $r.Hyperlinks.Add($r, $link.Addr)
    }
}

# --- View-state cosmetics (best effort) ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("F32").Select()
